$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table was missing the "Pr4NF" compound row, which belongs right
# after "Et4NF" (row 61) and before "Bu4NF" (old row 62). Insert a new row
# at position 62; this shifts all subsequent rows (old 62-113) down by one
# (to 63-114) and extends the used range from A1:I113 to A1:I114.
$ws.Rows.Item(62).Insert()

# Fill in the newly inserted row 62 with the Pr4NF parameters.
$ws.Range("A62").Value = "Pr4NF"
$ws.Range("B62").Value = 0.4463
$ws.Range("C62").Value = 0.409
$ws.Range("D62").Value = 0.0537
$ws.Range("E62").Value = 2
$ws.Range("F62").Value = 1
$ws.Range("G62").Value = -1
$ws.Range("H62").Value = 1
$ws.Range("I62").Value = 1

# Match the saved view/selection state.
$ws.Activate()
$aw = $excel.ActiveWindow
$aw.ScrollRow = 50
$aw.ScrollColumn = 1
$ws.Range("E62").Select()
